$wb = $excel.ActiveWorkbook

# Sheet "DB Test Result" holds the candidate test results in column D.
$ws1 = $wb.Worksheets.Item("DB Test Result")

$ws1.Range("D6").Value = 26
$ws1.Range("D7").Value = 26
$ws1.Range("D8").Value = 41
$ws1.Range("D9").Value = 43
$ws1.Range("D10").Value = 38
$ws1.Range("D11").Value = 39
$ws1.Range("D12").Value = 34

# Make "DB Test Result" the active sheet/tab and select F18 on it,
# which moves tabSelected from "DB Part 1 Project " to "DB Test Result".
$ws1.Activate()
$ws1.Range("F18").Select()
